# Update countries & provincias Spain
# Refresh COVID-19 country statistics table and fix sort order for a couple
# of countries whose totals crossed each other.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Update the "last refreshed" timestamp shown in A1.
$ws.Range("A1").Value = "Datos actualizados a 4 de Octubre de 2020 a las 21:29"

# --- Rows whose country stays the same, only the statistics refresh ---

# Row 4: Estados Unidos
Set-Row 4 7624081 22304 4839173 2570430 0 201 214478

# Row 5: India
Set-Row 5 6622135 74722 5583405 936016 0 902 102714

# Row 14: Francia
Set-Row 14 619190 12565 97778 489182 0 32 32230

# Row 26: Alemania
Set-Row 26 301013 985 261900 29512 0 4 9601

# Row 59: Uzbekistan
Set-Row 59 58612 374 55281 2851 0 3 480

# Row 94: Zambia
Set-Row 94 15052 78 14187 532 0 0 333

# Row 114: Zimbabue
Set-Row 114 7888 3 6359 1301 0 0 228

# Row 131: Ruanda
Set-Row 131 4866 14 3216 1621 0 0 29

# Row 137: Aruba
Set-Row 137 4079 5 3532 517 0 0 30

# Row 149: Sudan del Sur
Set-Row 149 2726 11 1290 1386 0 0 50

# --- Rows that swap order because their totals crossed each other ---

# Etiopia overtakes Honduras: row 54 becomes Etiopia, row 55 becomes Honduras.
$ws.Cells.Item(54, 1).Value = "Etiopia"
Set-Row 54 78819 959 33060 44537 0 8 1222

$ws.Cells.Item(55, 1).Value = "Honduras"
Set-Row 55 78788 519 29187 47202 0 13 2399

# Islas Malvinas overtakes Montserrat: row 215 becomes Islas Malvinas, row 216 becomes Montserrat.
$ws.Cells.Item(215, 1).Value = "Islas Malvinas"
Set-Row 215 13 0 13 0 0 0 0

$ws.Cells.Item(216, 1).Value = "Montserrat"
Set-Row 216 13 0 12 0 0 0 1
